$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clara")

# Refresh the two stale "wishlist image link" cells with new URLs (B10's
# link is untouched - only B13 and B9 get new images). Write B13 first so
# the shared-string table gets the new strings appended in the same order
# as the target file (B13's new string, then B9's new string).
$ws.Range("B13").Value = "https://www.galaxus.ch/im/productimages/4/8/5/0/4/6/4/9/2/5/9/8/2/2/9/9/7/1/9/0d7b50f8-1fb0-4c3c-9836-4f180c96c174_cropped.jpg"
$ws.Range("B9").Value = "https://m.media-amazon.com/images/I/71Ne9LAW7qL._AC_SL1500_.jpg"

# Match the author's final selection (B10) instead of the original E1.
$ws.Range("B10").Select()
